# Update system_info.qmd rendering to include code output.
#
# Before: paragraph 4 (FirstParagraph) has the two "import" lines;
#         paragraph 5 (BodyText) has the four print(...) lines;
#         paragraph 6 (BodyText) has the mem=... / print(...) lines.
# After:  paragraph 4 (SourceCode) holds ALL of the source code
#         (imports + prints + mem lines), syntax-highlighted with the
#         pandoc/highlighting character styles, using real line breaks;
#         paragraph 5 (SourceCode) holds the program's output, one line
#         per w:br, styled with VerbatimChar;
#         paragraph 6 is gone (its code merged into paragraph 4).

$d = $word.ActiveDocument

function Add-Run($text, $style) {
    $ins = $d.Range($global:pos, $global:pos)
    $ins.InsertAfter($text)
    $global:pos = $global:pos + $text.Length
    if ($style) {
        $sr = $d.Range($ins.Start, $global:pos)
        $sr.Style = $style
    }
}

function Add-Break() {
    $insB = $d.Range($global:pos, $global:pos)
    $insB.InsertBreak(6)   # wdLineBreak
    $global:pos = $global:pos + 1
}

# ------------------------------------------------------------------
# 1) Remove paragraph 6 ("mem = psutil.virtual_memory()" / BodyText)
#    entirely -- its code text gets folded into paragraph 4 below.
#    Deleted first (working from the end) so paragraph 4/5 indices
#    are unaffected.
# ------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Delete()

# ------------------------------------------------------------------
# 2) Paragraph 5 ("print(...)" lines / BodyText) -> becomes the
#    SourceCode/VerbatimChar *output* paragraph.
# ------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5clear = $d.Range($r5.Start, $r5.End - 1)
$r5clear.Delete()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Style = "SourceCode"

$global:pos = $d.Paragraphs.Item(5).Range.Start

Add-Run "Operating System: Linux 6.14.0-1011-aws" "VerbatimChar"
Add-Break
Add-Run "Python Version: 3.12.3" "VerbatimChar"
Add-Break
Add-Run "Machine: x86_64" "VerbatimChar"
Add-Break
Add-Run "Processor: x86_64" "VerbatimChar"
Add-Break
Add-Run "Total Memory (MB): 7938.12" "VerbatimChar"
Add-Break
Add-Run "Available Memory (MB): 4678.6" "VerbatimChar"

# ------------------------------------------------------------------
# 3) Paragraph 4 ("import platform" / "import psutil" / FirstParagraph)
#    -> becomes the SourceCode paragraph with ALL of the code,
#    syntax-highlighted and separated by real line breaks.
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4clear = $d.Range($r4.Start, $r4.End - 1)
$r4clear.Delete()
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Style = "SourceCode"

$global:pos = $d.Paragraphs.Item(4).Range.Start

Add-Run "import" "ImportTok"
Add-Run " platform" "NormalTok"
Add-Break
Add-Run "import" "ImportTok"
Add-Run " psutil" "NormalTok"
Add-Break
Add-Break
Add-Run "print" "BuiltInTok"
Add-Run "(" "NormalTok"
Add-Run '"Operating System:"' "StringTok"
Add-Run ", platform.system(), platform.release())" "NormalTok"
Add-Break
Add-Run "print" "BuiltInTok"
Add-Run "(" "NormalTok"
Add-Run '"Python Version:"' "StringTok"
Add-Run ", platform.python_version())" "NormalTok"
Add-Break
Add-Run "print" "BuiltInTok"
Add-Run "(" "NormalTok"
Add-Run '"Machine:"' "StringTok"
Add-Run ", platform.machine())" "NormalTok"
Add-Break
Add-Run "print" "BuiltInTok"
Add-Run "(" "NormalTok"
Add-Run '"Processor:"' "StringTok"
Add-Run ", platform.processor())" "NormalTok"
Add-Break
Add-Break
Add-Run "mem " "NormalTok"
Add-Run "=" "OperatorTok"
Add-Run " psutil.virtual_memory()" "NormalTok"
Add-Break
Add-Run "print" "BuiltInTok"
Add-Run "(" "NormalTok"
Add-Run '"Total Memory (MB):"' "StringTok"
Add-Run ", " "NormalTok"
Add-Run "round" "BuiltInTok"
Add-Run "(mem.total " "NormalTok"
Add-Run "/" "OperatorTok"
Add-Run " (" "NormalTok"
Add-Run "1024" "DecValTok"
Add-Run "**" "OperatorTok"
Add-Run "2" "DecValTok"
Add-Run "), " "NormalTok"
Add-Run "2" "DecValTok"
Add-Run "))" "NormalTok"
Add-Break
Add-Run "print" "BuiltInTok"
Add-Run "(" "NormalTok"
Add-Run '"Available Memory (MB):"' "StringTok"
Add-Run ", " "NormalTok"
Add-Run "round" "BuiltInTok"
Add-Run "(mem.available " "NormalTok"
Add-Run "/" "OperatorTok"
Add-Run " (" "NormalTok"
Add-Run "1024" "DecValTok"
Add-Run "**" "OperatorTok"
Add-Run "2" "DecValTok"
Add-Run "), " "NormalTok"
Add-Run "2" "DecValTok"
Add-Run "))" "NormalTok"

Write-Output "done"
